# Bill of Materials.xlsx - apply "Add files via upload" edit
#
# Summary of changes:
#  - Sheet "3D Printed Parts": fix a handful of typos in part names
#  - Sheet "Electronics": add a new BOM row for "Wires"
#  - Selection/active-cell bookkeeping per sheet, and the Electronics sheet
#    becomes the active tab (was Mechanical Hardware)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "3D Printed Parts" -- typo fixes
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("B3").Value  = "Motor & Arduino Holder [MIRRORED] (Only Use if v.1 is not used)"
$ws1.Range("B11").Value = "Centrifuge Test Tube Holder - No Holes"
$ws1.Range("B12").Value = "Centrifuge Test Tube Holder - With Holes"
$ws1.Range("B14").Value = "Support Piece for Digit Display and Button Board"

# ---------------------------------------------------------------------
# Sheet 2: "Electronics" -- new row for Wires
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A12").Value = 11
$ws2.Range("B12").Value = "Wires"
$ws2.Range("C12").Value = "Any electric wires that solves the connects properly (Various lengths)"
$ws2.Range("D12").Value = 23
$ws2.Range("D12").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# Selections / active sheet bookkeeping.
# Select in the order that leaves "Electronics" as the final active tab
# (selecting a range activates its sheet, and the last selection wins).
# ---------------------------------------------------------------------
$ws1.Range("C17").Select()

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("E18").Select()

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B27").Select()

$ws2.Range("C19").Select()
